$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New SKU values appended to the bottom of the list (rows 92-97)
$newSkus = @(10057381, 10092212, 10095739, 10020308, 10021623, 10029283)

$startRow = 92
for ($i = 0; $i -lt $newSkus.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newSkus[$i]
}

# Last row (97) reuses the same "highlighted" look already used on rows 89/91
# (Open Sans 9pt grey font, taller 15.75pt row).
$lastRow = $startRow + $newSkus.Length - 1
$ws.Cells.Item($lastRow, 1).Font.Name = "Open Sans"
$ws.Cells.Item($lastRow, 1).Font.Size = 9
$ws.Cells.Item($lastRow, 1).Font.Color = 4473924
$ws.Rows.Item($lastRow).RowHeight = 15.75

# Move the selection/view down to the newly added data, like the author did.
[void]$ws.Range("A97").Select()
